$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.931.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.648.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.46%  "
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0871"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.881.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.639.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.932.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.75%  "
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.453.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E37").Value = "  +3.35%  "
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("E42").Value = "  +2.47%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.790.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("E51").Value = "  +1.01%  "
